# invoices/10003-2023.1.18.xlsx
# Commit: "Add: tables of products information in pdf files"
#
# The underlying change re-lays-out the product table: the "total_price"
# header (E1) is switched from the numeric (#,##0.00) header format to the
# plain text header format already used by the other text headers
# (product_name / B1, and the blank F1 spacer cell), and every data row's
# height grows slightly (consistent with the new header format / column
# layout causing Excel to re-autofit the wrapped row heights).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights (autofit result of the relayout) ---------------------
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 33.75
$ws.Rows.Item(4).RowHeight = 19.5
$ws.Rows.Item(5).RowHeight = 19.5
$ws.Rows.Item(6).RowHeight = 19.5
$ws.Rows.Item(7).RowHeight = 19.5
$ws.Rows.Item(8).RowHeight = 19.5
$ws.Rows.Item(9).RowHeight = 20.25

# --- E1 "total_price" header: match the other text headers ------------
# Copy the format from B1 ("product_name" header) onto E1 so it stops
# using the numeric-header style and instead uses the same general/text
# header style as B1 and F1.
$ws.Range("B1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
